$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 199
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 199
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 199
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -425
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 17500
$ws.Range("J7").Value = 17500
$ws.Range("L7").Value = 17500
$ws.Range("N7").Value = -17724
$ws.Range("H14").Value = 17500
$ws.Range("J14").Value = 17500
$ws.Range("L14").Value = 17500
$ws.Range("N14").Value = -17882
$ws.Range("H21").Value = 70019
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -69551
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 70019
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -69785
$ws.Range("N23").ClearContents()
$ws.Range("H40").Value = 1576.3636
$ws.Range("I40").Value = 1532.3143
$ws.Range("J40").Value = 1747.6666
$ws.Range("K40").Value = 1532.3143
$ws.Range("L40").Value = 1747.6666
$ws.Range("M40").Value = -1357.3143
$ws.Range("N40").Value = -2097.6666
$ws.Range("H62").Value = 17859448
$ws.Range("I62").Value = 29413784
$ws.Range("J62").Value = 2747.4546
$ws.Range("K62").Value = 29413784
$ws.Range("L62").Value = 2747.4546
$ws.Range("M62").Value = -29413160
$ws.Range("N62").Value = -3995.4546
$ws.Range("H65").Value = 17859448
$ws.Range("I65").Value = 29413784
$ws.Range("J65").Value = 2747.4546
$ws.Range("K65").Value = 147068920
$ws.Range("L65").Value = 13737.273
$ws.Range("M65").Value = -147065800
$ws.Range("N65").Value = -19977.273
$ws.Range("H103").Value = 1533.9286
$ws.Range("I103").Value = 1033.3334
$ws.Range("J103").Value = 1670.4546
$ws.Range("K103").Value = 3100.0002
$ws.Range("L103").Value = 5011.3638
$ws.Range("M103").Value = -2514.0002
$ws.Range("N103").Value = -6183.3638
$ws.Range("H113").Value = 1752.4286
$ws.Range("I113").Value = 1655.3636
$ws.Range("J113").Value = 1796.9166
$ws.Range("K113").Value = 1655.3636
$ws.Range("L113").Value = 1796.9166
$ws.Range("M113").Value = 1598.6364
$ws.Range("N113").Value = -8304.9166
$ws.Range("H126").Value = 61490
$ws.Range("J126").Value = 61490
$ws.Range("L126").Value = 61490
$ws.Range("N126").Value = -71370
$ws.Range("H141").Value = 977.4231
$ws.Range("I141").Value = 713.6087
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 2140.8261
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 3039.1739
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 759999.75
$ws.Range("J6").Value = 13333
$ws.Range("L6").Value = 13333
$ws.Range("N6").Value = -13679
$ws.Range("H16").Value = 14000
$ws.Range("J16").Value = 14000
$ws.Range("L16").Value = 14000
$ws.Range("N16").Value = -14574
$ws.Range("H32").Value = 17249056
$ws.Range("I32").Value = 21278802
$ws.Range("J32").Value = 31054.455
$ws.Range("K32").Value = 21278802
$ws.Range("L32").Value = 31054.455
$ws.Range("M32").Value = -21278515
$ws.Range("N32").Value = -31628.455
$ws.Range("H74").Value = 9504.167
$ws.Range("I74").Value = 25759.25
$ws.Range("J74").Value = 1376.625
$ws.Range("K74").Value = 25759.25
$ws.Range("L74").Value = 1376.625
$ws.Range("M74").Value = -24885.25
$ws.Range("N74").Value = -3124.625
$ws.Range("H77").Value = 9504.167
$ws.Range("I77").Value = 25759.25
$ws.Range("J77").Value = 1376.625
$ws.Range("K77").Value = 128796.25
$ws.Range("L77").Value = 6883.125
$ws.Range("M77").Value = -124428.25
$ws.Range("N77").Value = -15619.125
$ws.Range("H114").Value = 23956.334
$ws.Range("J114").Value = 23956.334
$ws.Range("L114").Value = 23956.334
$ws.Range("N114").Value = -32634.334
$ws.Range("H119").Value = 28700
$ws.Range("J119").Value = 28700
$ws.Range("L119").Value = 28700
$ws.Range("N119").Value = -38376

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 28314
$ws.Range("J76").Value = 28314
$ws.Range("L76").Value = 28314
$ws.Range("N76").Value = -28944
$ws.Range("H79").Value = 28314
$ws.Range("J79").Value = 28314
$ws.Range("L79").Value = 28314
$ws.Range("N79").Value = -30498
$ws.Range("H105").Value = 2496.28
$ws.Range("I105").Value = 1201.5416
$ws.Range("J105").Value = 2905.1448
$ws.Range("K105").Value = 1201.5416
$ws.Range("L105").Value = 2905.1448
$ws.Range("M105").Value = 545.4584
$ws.Range("N105").Value = -6399.1448
$ws.Range("H107").Value = 4111.1113
$ws.Range("I107").Value = 3583.3333
$ws.Range("K107").Value = 3583.3333
$ws.Range("M107").Value = -1663.3333
$ws.Range("H134").Value = 4956
$ws.Range("I134").Value = 4956
$ws.Range("K134").Value = 14868
$ws.Range("M134").Value = -12333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 22000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 1846.4595
$ws.Range("I31").Value = 2036.25
$ws.Range("J31").Value = 1755.36
$ws.Range("K31").Value = 2036.25
$ws.Range("L31").Value = 1755.36
$ws.Range("M31").Value = -1741.25
$ws.Range("N31").Value = -2345.36
$ws.Range("H34").Value = 1846.4595
$ws.Range("I34").Value = 2036.25
$ws.Range("J34").Value = 1755.36
$ws.Range("K34").Value = 2036.25
$ws.Range("L34").Value = 1755.36
$ws.Range("M34").Value = -1834.25
$ws.Range("N34").Value = -2159.36
$ws.Range("H59").Value = 24500
$ws.Range("J59").Value = 24500
$ws.Range("L59").Value = 24500
$ws.Range("N59").Value = -26790
$ws.Range("H99").Value = 1725.7567
$ws.Range("I99").Value = 1352.238
$ws.Range("J99").Value = 2216
$ws.Range("K99").Value = 1352.238
$ws.Range("L99").Value = 2216
$ws.Range("M99").Value = 145.7619999999999
$ws.Range("N99").Value = -5212
$ws.Range("H126").Value = 1725.7567
$ws.Range("I126").Value = 1352.238
$ws.Range("J126").Value = 2216
$ws.Range("K126").Value = 4056.714
$ws.Range("L126").Value = 6648
$ws.Range("M126").Value = -1586.714
$ws.Range("N126").Value = -11588
$ws.Range("H132").Value = 1795.8846
$ws.Range("I132").Value = 918.35
$ws.Range("K132").Value = 2755.05
$ws.Range("M132").Value = -225.0500000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -3224
$ws.Range("H6").Value = 312.69232
$ws.Range("I6").Value = 312.69232
$ws.Range("K6").Value = 938.07696
$ws.Range("M6").Value = -825.07696
$ws.Range("H122").Value = 8490.667
$ws.Range("I122").Value = 16758.77
$ws.Range("J122").Value = 813.1429
$ws.Range("K122").Value = 150828.93
$ws.Range("L122").Value = 7318.2861
$ws.Range("M122").Value = -148378.93
$ws.Range("N122").Value = -12218.2861
$ws.Range("H127").Value = 733.3333
$ws.Range("J127").Value = 733.3333
$ws.Range("L127").Value = 2199.9999
$ws.Range("N127").Value = -12119.9999
$ws.Range("H135").Value = 1000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 9000
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8971
$ws.Range("I5").Value = 3574.6667
$ws.Range("J5").Value = 11669.167
$ws.Range("K5").Value = 3574.6667
$ws.Range("L5").Value = 11669.167
$ws.Range("M5").Value = -3462.6667
$ws.Range("N5").Value = -11893.167
$ws.Range("H126").Value = 1247
$ws.Range("I126").Value = 1083
$ws.Range("J126").Value = 1575
$ws.Range("K126").Value = 3249
$ws.Range("L126").Value = 4725
$ws.Range("M126").Value = -779
$ws.Range("N126").Value = -9665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 48239.5
$ws.Range("J119").Value = 48239.5
$ws.Range("L119").Value = 48239.5
$ws.Range("N119").Value = -57915.5
$ws.Range("H122").Value = 2437.6956
$ws.Range("I122").Value = 1926.8823
$ws.Range("J122").Value = 3885
$ws.Range("K122").Value = 5780.6469
$ws.Range("L122").Value = 11655
$ws.Range("M122").Value = -3330.6469
$ws.Range("N122").Value = -16555
$ws.Range("H138").Value = 48000
$ws.Range("J138").Value = 48000
$ws.Range("L138").Value = 48000
$ws.Range("N138").Value = -58280
